$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.005243333333333
$ws.Range("H2").Value = 3.01573
$ws.Range("I2").Value = 0.07224874268505826
$ws.Range("J2").Value = 0.07224874268505825
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.438907666666667
$ws.Range("N2").Value = 10.316723
$ws.Range("O2").Value = 0.05825422340060618
$ws.Range("P2").Value = 0.05825422340060618
$ws.Range("Q2").Value = 3.456939005865555
$ws.Range("R2").Value = 31.11245105278999
$ws.Range("S2").Value = 0.004208794396788296
$ws.Range("T2").Value = 0.004208794396788295
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.005243333333333
$ws.Range("H3").Value = 3.01573
$ws.Range("I3").Value = 0.07224874268505826
$ws.Range("J3").Value = 0.07224874268505825
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 10.383857
$ws.Range("N3").Value = 31.151571
$ws.Range("O3").Value = 0.1758999031294962
$ws.Range("P3").Value = 0.1758999031294962
$ws.Range("Q3").Value = 10.43830302353667
$ws.Range("R3").Value = 93.94472721182998
$ws.Range("S3").Value = 0.01270854683952964
$ws.Range("T3").Value = 0.01270854683952964
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.005243333333333
$ws.Range("H4").Value = 3.01573
$ws.Range("I4").Value = 0.07224874268505826
$ws.Range("J4").Value = 0.07224874268505825
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 25.34077833333333
$ws.Range("N4").Value = 76.022335
$ws.Range("O4").Value = 0.4292663558501786
$ws.Range("P4").Value = 0.4292663558501786
$ws.Range("Q4").Value = 25.47364848106111
$ws.Range("R4").Value = 229.26283632955
$ws.Range("S4").Value = 0.03101395448717221
$ws.Range("T4").Value = 0.03101395448717221
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.005243333333333
$ws.Range("H5").Value = 3.01573
$ws.Range("I5").Value = 0.07224874268505826
$ws.Range("J5").Value = 0.07224874268505825
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 19.86921833333334
$ws.Range("N5").Value = 59.60765500000001
$ws.Range("O5").Value = 0.336579517619719
$ws.Range("P5").Value = 0.336579517619719
$ws.Range("Q5").Value = 19.97339926812778
$ws.Range("R5").Value = 179.76059341315
$ws.Range("S5").Value = 0.02431744696156811
$ws.Range("T5").Value = 0.02431744696156811
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 10.25983933333333
$ws.Range("H6").Value = 30.779518
$ws.Range("I6").Value = 0.7373940889775011
$ws.Range("J6").Value = 0.737394088977501
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.438907666666667
$ws.Range("N6").Value = 10.316723
$ws.Range("O6").Value = 0.05825422340060618
$ws.Range("P6").Value = 0.05825422340060618
$ws.Range("Q6").Value = 35.28264014216823
$ws.Range("R6").Value = 317.543761279514
$ws.Range("S6").Value = 0.04295631999358182
$ws.Range("T6").Value = 0.04295631999358181
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 10.25983933333333
$ws.Range("H7").Value = 30.779518
$ws.Range("I7").Value = 0.7373940889775011
$ws.Range("J7").Value = 0.737394088977501
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.383857
$ws.Range("N7").Value = 31.151571
$ws.Range("O7").Value = 0.1758999031294962
$ws.Range("P7").Value = 0.1758999031294962
$ws.Range("Q7").Value = 106.5367044803087
$ws.Range("R7").Value = 958.8303403227781
$ws.Range("S7").Value = 0.1297075488194055
$ws.Range("T7").Value = 0.1297075488194055
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 10.25983933333333
$ws.Range("H8").Value = 30.779518
$ws.Range("I8").Value = 0.7373940889775011
$ws.Range("J8").Value = 0.737394088977501
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 25.34077833333333
$ws.Range("N8").Value = 76.022335
$ws.Range("O8").Value = 0.4292663558501786
$ws.Range("P8").Value = 0.4292663558501786
$ws.Range("Q8").Value = 259.9923142816144
$ws.Range("R8").Value = 2339.93082853453
$ws.Range("S8").Value = 0.3165384734008342
$ws.Range("T8").Value = 0.3165384734008342
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 10.25983933333333
$ws.Range("H9").Value = 30.779518
$ws.Range("I9").Value = 0.7373940889775011
$ws.Range("J9").Value = 0.737394088977501
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.86921833333334
$ws.Range("N9").Value = 59.60765500000001
$ws.Range("O9").Value = 0.336579517619719
$ws.Range("P9").Value = 0.336579517619719
$ws.Range("Q9").Value = 203.8549877789212
$ws.Range("R9").Value = 1834.69489001029
$ws.Range("S9").Value = 0.2481917467636795
$ws.Range("T9").Value = 0.2481917467636795
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.383875
$ws.Range("H10").Value = 1.151625
$ws.Range("I10").Value = 0.02758982345723265
$ws.Range("J10").Value = 0.02758982345723265
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.438907666666667
$ws.Range("N10").Value = 10.316723
$ws.Range("O10").Value = 0.05825422340060618
$ws.Range("P10").Value = 0.05825422340060618
$ws.Range("Q10").Value = 1.320110680541667
$ws.Range("R10").Value = 11.880996124875
$ws.Range("S10").Value = 0.001607223739260916
$ws.Range("T10").Value = 0.001607223739260916
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.383875
$ws.Range("H11").Value = 1.151625
$ws.Range("I11").Value = 0.02758982345723265
$ws.Range("J11").Value = 0.02758982345723265
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 10.383857
$ws.Range("N11").Value = 31.151571
$ws.Range("O11").Value = 0.1758999031294962
$ws.Range("P11").Value = 0.1758999031294962
$ws.Range("Q11").Value = 3.986103105875001
$ws.Range("R11").Value = 35.87492795287501
$ws.Range("S11").Value = 0.004853047273487125
$ws.Range("T11").Value = 0.004853047273487125
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.383875
$ws.Range("H12").Value = 1.151625
$ws.Range("I12").Value = 0.02758982345723265
$ws.Range("J12").Value = 0.02758982345723265
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 25.34077833333333
$ws.Range("N12").Value = 76.022335
$ws.Range("O12").Value = 0.4292663558501786
$ws.Range("P12").Value = 0.4292663558501786
$ws.Range("Q12").Value = 9.727691282708333
$ws.Range("R12").Value = 87.549221544375
$ws.Range("S12").Value = 0.01184338297403604
$ws.Range("T12").Value = 0.01184338297403604
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.383875
$ws.Range("H13").Value = 1.151625
$ws.Range("I13").Value = 0.02758982345723265
$ws.Range("J13").Value = 0.02758982345723265
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 19.86921833333334
$ws.Range("N13").Value = 59.60765500000001
$ws.Range("O13").Value = 0.336579517619719
$ws.Range("P13").Value = 0.336579517619719
$ws.Range("Q13").Value = 7.627296187708335
$ws.Range("R13").Value = 68.64566568937502
$ws.Range("S13").Value = 0.009286169470448575
$ws.Range("T13").Value = 0.009286169470448576
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.264687
$ws.Range("H14").Value = 6.794061
$ws.Range("I14").Value = 0.162767344880208
$ws.Range("J14").Value = 0.162767344880208
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.438907666666667
$ws.Range("N14").Value = 10.316723
$ws.Range("O14").Value = 0.05825422340060618
$ws.Range("P14").Value = 0.05825422340060618
$ws.Range("Q14").Value = 7.788049486900333
$ws.Range("R14").Value = 70.09244538210299
$ws.Range("S14").Value = 0.009481885270975149
$ws.Range("T14").Value = 0.009481885270975149
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.264687
$ws.Range("H15").Value = 6.794061
$ws.Range("I15").Value = 0.162767344880208
$ws.Range("J15").Value = 0.162767344880208
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 10.383857
$ws.Range("N15").Value = 31.151571
$ws.Range("O15").Value = 0.1758999031294962
$ws.Range("P15").Value = 0.1758999031294962
$ws.Range("Q15").Value = 23.516185957759
$ws.Range("R15").Value = 211.645673619831
$ws.Range("S15").Value = 0.02863076019707388
$ws.Range("T15").Value = 0.02863076019707388
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.264687
$ws.Range("H16").Value = 6.794061
$ws.Range("I16").Value = 0.162767344880208
$ws.Range("J16").Value = 0.162767344880208
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 25.34077833333333
$ws.Range("N16").Value = 76.022335
$ws.Range("O16").Value = 0.4292663558501786
$ws.Range("P16").Value = 0.4292663558501786
$ws.Range("Q16").Value = 57.38893126138166
$ws.Range("R16").Value = 516.5003813524349
$ws.Range("S16").Value = 0.06987054498813611
$ws.Range("T16").Value = 0.06987054498813612
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2.264687
$ws.Range("H17").Value = 6.794061
$ws.Range("I17").Value = 0.162767344880208
$ws.Range("J17").Value = 0.162767344880208
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 19.86921833333334
$ws.Range("N17").Value = 59.60765500000001
$ws.Range("O17").Value = 0.336579517619719
$ws.Range("P17").Value = 0.336579517619719
$ws.Range("Q17").Value = 44.99756045966167
$ws.Range("R17").Value = 404.9780441369551
$ws.Range("S17").Value = 0.05478415442402285
$ws.Range("T17").Value = 0.05478415442402285
